# Natmi following Dr Hou advice
# The LR-pairs result table is recomputed: the sheet grows from 9 rows (8 data
# rows) to 13 rows (12 data rows) covering all Sending/Target cluster
# combinations (ECs, FAPs, M2, sCs) for the Gas6 -> Tyro3 ligand/receptor
# pair, and all numeric metrics are refreshed with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Sending cluster, Ligand, Receptor, Target cluster, then the
# 16 numeric metric columns (E..T).
$rows = @(
    @(2,  "ECs",  "Gas6", "Tyro3", "ECs",  3,1,15.79785166666667,47.39355500000001,0.1445757693628457,0.1445757693628457,1,0.3333333333333333,0.096929,0.290787,0.01940247578118474,0.01940247578118474,1.531269964198334,13.781429677785,0.002805127863608765,0.002805127863608765),
    @(3,  "ECs",  "Gas6", "Tyro3", "FAPs", 3,1,15.79785166666667,47.39355500000001,0.1445757693628457,0.1445757693628457,3,1,3.954952666666667,11.864858,0.7916709481242146,0.7916709481242146,62.47975557668779,562.31780019019,0.1144564364072718,0.1144564364072719),
    @(4,  "ECs",  "Gas6", "Tyro3", "sCs",  3,1,15.79785166666667,47.39355500000001,0.1445757693628457,0.1445757693628457,3,1,0.9438210000000001,2.831463,0.1889265760946008,0.1889265760946008,14.910344157885,134.193097420965,0.02731420509196513,0.02731420509196513),
    @(5,  "FAPs", "Gas6", "Tyro3", "ECs",  3,1,74.39645633333333,223.189369,0.6808473164079603,0.6808473164079603,1,0.3333333333333333,0.096929,0.290787,0.01940247578118474,0.01940247578118474,7.211174115933667,64.900567043403,0.01321012356729007,0.01321012356729007),
    @(6,  "FAPs", "Gas6", "Tyro3", "FAPs", 3,1,74.39645633333333,223.189369,0.6808473164079603,0.6808473164079603,3,1,3.954952666666667,11.864858,0.7916709481242146,0.7916709481242146,294.2344633660669,2648.110170294602,0.539007040508517,0.539007040508517),
    @(7,  "FAPs", "Gas6", "Tyro3", "sCs",  3,1,74.39645633333333,223.189369,0.6808473164079603,0.6808473164079603,3,1,0.9438210000000001,2.831463,0.1889265760946008,0.1889265760946008,70.21693781298301,631.9524403168471,0.1286301523321532,0.1286301523321532),
    @(8,  "M2",   "Gas6", "Tyro3", "ECs",  3,1,16.36992,49.10976,0.1498111153557632,0.1498111153557632,1,0.3333333333333333,0.096929,0.290787,0.01940247578118474,0.01940247578118474,1.58671997568,14.28047978112,0.002906706537442467,0.002906706537442467),
    @(9,  "M2",   "Gas6", "Tyro3", "FAPs", 3,1,16.36992,49.10976,0.1498111153557632,0.1498111153557632,3,1,3.954952666666667,11.864858,0.7916709481242146,0.7916709481242146,64.74225875712,582.68032881408,0.1186011077332431,0.1186011077332431),
    @(10, "M2",   "Gas6", "Tyro3", "sCs",  3,1,16.36992,49.10976,0.1498111153557632,0.1498111153557632,3,1,0.9438210000000001,2.831463,0.1889265760946008,0.1889265760946008,15.45027426432,139.05246837888,0.02830330108507761,0.02830330108507761),
    @(11, "sCs",  "Gas6", "Tyro3", "ECs",  3,1,2.706168666666667,8.118506,0.02476579887343077,0.02476579887343077,1,0.3333333333333333,0.096929,0.290787,0.01940247578118474,0.01940247578118474,0.2623062226913334,2.360756004222,0.0004805178128434327,0.0004805178128434327),
    @(12, "sCs",  "Gas6", "Tyro3", "FAPs", 3,1,2.706168666666667,8.118506,0.02476579887343077,0.02476579887343077,3,1,3.954952666666667,11.864858,0.7916709481242146,0.7916709481242146,10.70276898468311,96.32492086214801,0.01960636347518254,0.01960636347518254),
    @(13, "sCs",  "Gas6", "Tyro3", "sCs",  3,1,2.706168666666667,8.118506,0.02476579887343077,0.02476579887343077,3,1,0.9438210000000001,2.831463,0.1889265760946008,0.1889265760946008,2.554138817142,22.987249354278,0.004678917585404797,0.004678917585404797)
)

foreach ($row in $rows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]

    for ($i = 5; $i -lt $row.Length; $i++) {
        # $row[$i] for $i=5 holds column E (the 5th spreadsheet column),
        # so the column number equals the array index directly.
        $col = $i
        $ws.Cells.Item($r, $col).Value = $row[$i]
    }
}

Write-Output "Updated sheet to A1:T13 with refreshed Natmi Gas6-Tyro3 values"
